$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 61.526261
$ws.Range("H2").Value = 184.578783
$ws.Range("I2").Value = 0.9684015388399598
$ws.Range("J2").Value = 0.9684015388399598
$ws.Range("M2").Value = 3.704480666666667
$ws.Range("N2").Value = 11.113442
$ws.Range("O2").Value = 0.298964201395561
$ws.Range("P2").Value = 0.2989642013955609
$ws.Range("Q2").Value = 227.9228443667873
$ws.Range("R2").Value = 2051.305599301086
$ws.Range("S2").Value = 0.2895173926895209
$ws.Range("T2").Value = 0.2895173926895209

# Row 3
$ws.Range("G3").Value = 61.526261
$ws.Range("H3").Value = 184.578783
$ws.Range("I3").Value = 0.9684015388399598
$ws.Range("J3").Value = 0.9684015388399598
$ws.Range("O3").Value = 0.3765624616238499
$ws.Range("P3").Value = 0.3765624616238499
$ws.Range("Q3").Value = 287.0818209485514
$ws.Range("R3").Value = 2583.736388536962
$ws.Range("S3").Value = 0.3646636673058996
$ws.Range("T3").Value = 0.3646636673058996

# Row 4
$ws.Range("G4").Value = 61.526261
$ws.Range("H4").Value = 184.578783
$ws.Range("I4").Value = 0.9684015388399598
$ws.Range("J4").Value = 0.9684015388399598
$ws.Range("M4").Value = 2.870093333333334
$ws.Range("N4").Value = 8.610280000000001
$ws.Range("O4").Value = 0.2316263029934534
$ws.Range("P4").Value = 0.2316263029934534
$ws.Range("Q4").Value = 176.5861115210267
$ws.Range("R4").Value = 1589.27500368924
$ws.Range("S4").Value = 0.2243072682546711
$ws.Range("T4").Value = 0.2243072682546711

# Row 5
$ws.Range("G5").Value = 61.526261
$ws.Range("H5").Value = 184.578783
$ws.Range("I5").Value = 0.9684015388399598
$ws.Range("J5").Value = 0.9684015388399598
$ws.Range("M5").Value = 1.150472333333333
$ws.Range("N5").Value = 3.451417
$ws.Range("O5").Value = 0.09284703398713583
$ws.Range("P5").Value = 0.09284703398713583
$ws.Range("Q5").Value = 70.78426105394567
$ws.Range("R5").Value = 637.058349485511
$ws.Range("S5").Value = 0.08991321058986838
$ws.Range("T5").Value = 0.08991321058986838

# Row 6
$ws.Range("I6").Value = 0.002979850677668077
$ws.Range("J6").Value = 0.002979850677668078
$ws.Range("M6").Value = 3.704480666666667
$ws.Range("N6").Value = 11.113442
$ws.Range("O6").Value = 0.298964201395561
$ws.Range("P6").Value = 0.2989642013955609
$ws.Range("Q6").Value = 0.7013372191208889
$ws.Range("R6").Value = 6.312034972088
$ws.Range("S6").Value = 0.0008908686781270578
$ws.Range("T6").Value = 0.0008908686781270577

# Row 7
$ws.Range("I7").Value = 0.002979850677668077
$ws.Range("J7").Value = 0.002979850677668078
$ws.Range("O7").Value = 0.3765624616238499
$ws.Range("P7").Value = 0.3765624616238499
$ws.Range("S7").Value = 0.001122099906454189
$ws.Range("T7").Value = 0.001122099906454189

# Row 8
$ws.Range("I8").Value = 0.002979850677668077
$ws.Range("J8").Value = 0.002979850677668078
$ws.Range("M8").Value = 2.870093333333334
$ws.Range("N8").Value = 8.610280000000001
$ws.Range("O8").Value = 0.2316263029934534
$ws.Range("P8").Value = 0.2316263029934534
$ws.Range("Q8").Value = 0.5433698966577779
$ws.Range("R8").Value = 4.890329069920001
$ws.Range("S8").Value = 0.0006902117959407935
$ws.Range("T8").Value = 0.0006902117959407936

# Row 9
$ws.Range("I9").Value = 0.002979850677668077
$ws.Range("J9").Value = 0.002979850677668078
$ws.Range("M9").Value = 1.150472333333333
$ws.Range("N9").Value = 3.451417
$ws.Range("O9").Value = 0.09284703398713583
$ws.Range("P9").Value = 0.09284703398713583
$ws.Range("Q9").Value = 0.2178089561097778
$ws.Range("R9").Value = 1.960280604988
$ws.Range("S9").Value = 0.0002766702971460377
$ws.Range("T9").Value = 0.0002766702971460377

# Row 10
$ws.Range("G10").Value = 0.6472316666666668
$ws.Range("H10").Value = 1.941695
$ws.Range("I10").Value = 0.01018719700821657
$ws.Range("J10").Value = 0.01018719700821657
$ws.Range("M10").Value = 3.704480666666667
$ws.Range("N10").Value = 11.113442
$ws.Range("O10").Value = 0.298964201395561
$ws.Range("P10").Value = 0.2989642013955609
$ws.Range("Q10").Value = 2.397657196021111
$ws.Range("R10").Value = 21.57891476419
$ws.Range("S10").Value = 0.003045607218020716
$ws.Range("T10").Value = 0.003045607218020715

# Row 11
$ws.Range("G11").Value = 0.6472316666666668
$ws.Range("H11").Value = 1.941695
$ws.Range("I11").Value = 0.01018719700821657
$ws.Range("J11").Value = 0.01018719700821657
$ws.Range("O11").Value = 0.3765624616238499
$ws.Range("P11").Value = 0.3765624616238499
$ws.Range("Q11").Value = 3.019985977081112
$ws.Range("R11").Value = 27.17987379373
$ws.Range("S11").Value = 0.003836115982461152
$ws.Range("T11").Value = 0.003836115982461152

# Row 12
$ws.Range("G12").Value = 0.6472316666666668
$ws.Range("H12").Value = 1.941695
$ws.Range("I12").Value = 0.01018719700821657
$ws.Range("J12").Value = 0.01018719700821657
$ws.Range("M12").Value = 2.870093333333334
$ws.Range("N12").Value = 8.610280000000001
$ws.Range("O12").Value = 0.2316263029934534
$ws.Range("P12").Value = 0.2316263029934534
$ws.Range("Q12").Value = 1.857615291622223
$ws.Range("R12").Value = 16.7185376246
$ws.Range("S12").Value = 0.002359622780879174
$ws.Range("T12").Value = 0.002359622780879174

# Row 13
$ws.Range("G13").Value = 0.6472316666666668
$ws.Range("H13").Value = 1.941695
$ws.Range("I13").Value = 0.01018719700821657
$ws.Range("J13").Value = 0.01018719700821657
$ws.Range("M13").Value = 1.150472333333333
$ws.Range("N13").Value = 3.451417
$ws.Range("O13").Value = 0.09284703398713583
$ws.Range("P13").Value = 0.09284703398713583
$ws.Range("Q13").Value = 0.7446221257572223
$ws.Range("R13").Value = 6.701599131815001
$ws.Range("S13").Value = 0.0009458510268555326
$ws.Range("T13").Value = 0.0009458510268555326

# Row 14
$ws.Range("G14").Value = 1.171018333333333
$ws.Range("H14").Value = 3.513055
$ws.Range("I14").Value = 0.01843141347415545
$ws.Range("J14").Value = 0.01843141347415545
$ws.Range("M14").Value = 3.704480666666667
$ws.Range("N14").Value = 11.113442
$ws.Range("O14").Value = 0.298964201395561
$ws.Range("P14").Value = 0.2989642013955609
$ws.Range("Q14").Value = 4.338014776145555
$ws.Range("R14").Value = 39.04213298531
$ws.Range("S14").Value = 0.005510332809892267
$ws.Range("T14").Value = 0.005510332809892265

# Row 15
$ws.Range("G15").Value = 1.171018333333333
$ws.Range("H15").Value = 3.513055
$ws.Range("I15").Value = 0.01843141347415545
$ws.Range("J15").Value = 0.01843141347415545
$ws.Range("O15").Value = 0.3765624616238499
$ws.Range("P15").Value = 0.3765624616238499
$ws.Range("Q15").Value = 5.463977008085555
$ws.Range("R15").Value = 49.17579307277001
$ws.Range("S15").Value = 0.006940578429034973
$ws.Range("T15").Value = 0.006940578429034973

# Row 16
$ws.Range("G16").Value = 1.171018333333333
$ws.Range("H16").Value = 3.513055
$ws.Range("I16").Value = 0.01843141347415545
$ws.Range("J16").Value = 0.01843141347415545
$ws.Range("M16").Value = 2.870093333333334
$ws.Range("N16").Value = 8.610280000000001
$ws.Range("O16").Value = 0.2316263029934534
$ws.Range("P16").Value = 0.2316263029934534
$ws.Range("Q16").Value = 3.360931911711111
$ws.Range("R16").Value = 30.2483872054
$ws.Range("S16").Value = 0.004269200161962351
$ws.Range("T16").Value = 0.004269200161962351

# Row 17
$ws.Range("G17").Value = 1.171018333333333
$ws.Range("H17").Value = 3.513055
$ws.Range("I17").Value = 0.01843141347415545
$ws.Range("J17").Value = 0.01843141347415545
$ws.Range("M17").Value = 1.150472333333333
$ws.Range("N17").Value = 3.451417
$ws.Range("O17").Value = 0.09284703398713583
$ws.Range("P17").Value = 0.09284703398713583
$ws.Range("Q17").Value = 1.347224194326111
$ws.Range("R17").Value = 12.125017748935
$ws.Range("S17").Value = 0.001711302073265864
$ws.Range("T17").Value = 0.001711302073265864
